$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.604.19"
$ws.Range("E2").Value = "  +0.17%  "
$ws.Range("D3").Value = "3.522.52"
$ws.Range("E3").Value = "  -1.97%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "606.91"
$ws.Range("E5").Value = "  -0.30%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.26"
$ws.Range("E6").Value = "  -3.91%  "
$ws.Range("D7").Value = "3.521.32"
$ws.Range("E7").Value = "  -2.01%  "
$ws.Range("E8").Value = "  -0.11%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.510"
$ws.Range("E9").Value = "  +4.08%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.73"
$ws.Range("E10").Value = "  -3.55%  "
$ws.Range("E11").Value = "  -4.47%  "
$ws.Range("E12").Value = "  -1.73%  "
$ws.Range("D13").Value = "4.111.72"
$ws.Range("E13").Value = "  -2.20%  "
$ws.Range("E14").Value = "  -6.42%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "28.68"
$ws.Range("E15").Value = "  -3.81%  "
$ws.Range("D16").Value = "3.521.74"
$ws.Range("E16").Value = "  -2.00%  "
$ws.Range("E17").Value = "  +0.37%  "
$ws.Range("D18").Value = "66.443.04"
$ws.Range("E18").Value = "  -0.18%  "
$ws.Range("E19").Value = "  -6.92%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.14"
$ws.Range("E20").Value = "  -3.31%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.62"
$ws.Range("E21").Value = "  -3.16%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "423.02"
$ws.Range("E22").Value = "  -1.00%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.588"
$ws.Range("E23").Value = "  -5.07%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "77.03"
$ws.Range("E24").Value = "  -2.12%  "
$ws.Range("D25").Value = "3.663.43"
$ws.Range("E25").Value = "  -2.06%  "
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("E27").Value = "  -5.24%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.92"
$ws.Range("E28").Value = "  -4.74%  "
$ws.Range("E29").Value = "  -2.02%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.95"
$ws.Range("E30").Value = "  -4.91%  "
$ws.Range("E31").Value = "  +0.25%  "
$ws.Range("D32").Value = "3.526.56"
$ws.Range("E32").Value = "  -1.78%  "
$ws.Range("E33").Value = "  -2.24%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "24.16"
$ws.Range("E34").Value = "  -4.98%  "
$ws.Range("E36").Value = "  -9.59%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "7.58"
$ws.Range("E37").Value = "  -3.53%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.62"
$ws.Range("E38").Value = "  -4.05%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "173.62"
$ws.Range("E39").Value = "  -2.15%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.21"
$ws.Range("E40").Value = "  -7.58%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0813"
$ws.Range("E41").Value = "  -4.95%  "
$ws.Range("E42").Value = "  -4.56%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.851"
$ws.Range("E43").Value = "  -5.08%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "45.49"
$ws.Range("E44").Value = "  -0.85%  "
$ws.Range("E45").Value = "  -6.59%  "
$ws.Range("E46").Value = "  +0.05%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.36"
$ws.Range("E47").Value = "  -8.05%  "
$ws.Range("E48").Value = "  -1.74%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.12"
$ws.Range("E49").Value = "  -4.35%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "22.88"
$ws.Range("E50").Value = "  -4.74%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.904"
$ws.Range("E51").Value = "  -5.09%  "
